$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.580.59'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.884.85'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.52'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.71%  '

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4738'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2895'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -1.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06540'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.34'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +1.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7762'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +5.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '101.06'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +4.32%  '

$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.883.61'
$ws.Range("D14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.261'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.33'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +0.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.556.38'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -0.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.22'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007531'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -0.19%  '

$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.128.91'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.351'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +0.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9999'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.429'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +2.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.165'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -0.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.31'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -1.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.13'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.916'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.18%  '

$ws.Range("E29").Value = '  -0.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09705'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -0.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.504'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +0.83%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.264'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -0.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.195'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("E34").Value = '  -0.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.131'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +0.31%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6974'
$ws.Range("D36").NumberFormat = "General"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.756'
$ws.Range("D37").NumberFormat = "General"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01911'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +0.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.897'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +3.23%  '

$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.304'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.984'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -0.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4262'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -0.46%  '

$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8312'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -0.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.52'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -0.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.852'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +3.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.030'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("E49").Value = '  -1.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '893.23'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -2.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05778'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +0.42%  '
